# Stereo record implemented from Camera config
#
# Flip the "video_file" (B3) and "vid_pose_file" (B5) flags on the
# Test_1 sheet from 0 to 1, enabling stereo recording driven by the
# camera configuration. Also re-activate Test_1 (it is the tab the user
# was on) and move the selection to B5 to match where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_1")

# Keep Test_1 as the active/visible sheet (it was tabSelected in the
# source workbook).
$ws.Activate()

$ws.Range("B3").Value = 1
$ws.Range("B5").Value = 1

# Leave the selection on B5, matching the post-edit cursor position.
$ws.Range("B5").Select()
